$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns being updated: B, C, D, E, F, I, J, K, L, M, N (G stays 1, H is blank, A is index)
$cols = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

$data = @(
    @(2, @(1.02, 1.085285357373827, 1.075804248105593, 1.085790653579353, 1.091945119229056, 1.054874438299318, 1.090140975215735, 1.078489851521646, 1.088450171361276, 1.094588817171968, 1.091689100441459)),
    @(3, @(1.02, 1.087727981886532, 1.077035246297037, 1.087832406976849, 1.093713350435646, 1.055339924768057, 1.092239953658483, 1.079537221343549, 1.090308249021628, 1.096175191546702, 1.093791059674352)),
    @(4, @(1.02, 1.089302210823541, 1.077827588592081, 1.089147715689942, 1.094852229045761, 1.055637481483805, 1.09359174519508, 1.080210182827273, 1.091504310248888, 1.09719593976768, 1.095144770909984)),
    @(5, @(1.02, 1.089962545482188, 1.078159700573001, 1.089699305117754, 1.095329774034578, 1.055761709754973, 1.094158543630963, 1.080491972232219, 1.092005671303806, 1.097623712110326, 1.095712374264658)),
    @(6, @(1.02, 1.090073333461648, 1.078215406111953, 1.089791840218525, 1.095409883992166, 1.055782517755608, 1.094253624928233, 1.080539220462771, 1.092089767118702, 1.097695458497801, 1.095807590588276)),
    @(7, @(1.02, 1.089311039974704, 1.077832030152311, 1.089155091380388, 1.094858614861395, 1.055639144814616, 1.093599324610461, 1.08021395250559, 1.091511015167703, 1.097201660955182, 1.095152361089003)),
    @(8, @(1.02, 1.086112185062391, 1.076221148058232, 1.086481902953986, 1.092543811175247, 1.055032510940251, 1.09085167839765, 1.078844809207014, 1.089079427840996, 1.095126143584516, 1.092400812903411)),
    @(9, @(1.02, 1.080425233092606, 1.073349715429244, 1.08172519410482, 1.088423185686085, 1.053935277475849, 1.08595948464054, 1.07639508135159, 1.084745486745364, 1.091423729042012, 1.087501671670023)),
    @(10, @(1.02, 1.076597644989844, 1.071412322333998, 1.078520918950357, 1.08564645213221, 1.053184276542904, 1.08266186395706, 1.074735991687704, 1.08182117119621, 1.088923568601901, 1.084199367986931)),
    @(11, @(1.02, 1.07493107349347, 1.07056769949168, 1.077125107267864, 1.084436691011298, 1.052854341936839, 1.081224883321466, 1.07401121390971, 1.080546162000244, 1.087833054761014, 1.082760346674011)),
    @(12, @(1.02, 1.074310602245474, 1.070253089660974, 1.076605347276907, 1.083986185106333, 1.052731066388176, 1.080689715889926, 1.073741021136204, 1.080071212132816, 1.087426767430806, 1.082224419243333)),
    @(13, @(1.02, 1.074443761068728, 1.070320614655497, 1.076716896746665, 1.084082872573371, 1.052757542299945, 1.080804575578681, 1.073799022977018, 1.080173152418986, 1.087513973240008, 1.08233944204602)),
    @(14, @(1.02, 1.074879814660227, 1.070541711800221, 1.077082170366509, 1.084399475636039, 1.052844166752545, 1.081180675208115, 1.073988899755865, 1.080506930355972, 1.087799496031397, 1.082716075780067)),
    @(15, @(1.02, 1.07514829034927, 1.07067781992406, 1.077307054889626, 1.08459439229524, 1.052897442808034, 1.081412214608342, 1.074105758921001, 1.080712401289627, 1.08797525319522, 1.0829479439928)),
    @(16, @(1.02, 1.076708051749913, 1.071468254973639, 1.07861337523482, 1.085726580969634, 1.053206072362823, 1.082757036454049, 1.074783956717545, 1.081905601506107, 1.088995772762099, 1.084294675639782)),
    @(17, @(1.02, 1.077683950528012, 1.071962528829719, 1.079430534436532, 1.086434764535959, 1.053398389775249, 1.083598143933855, 1.075207650353503, 1.082651691951928, 1.089633773605885, 1.085136977588685)),
    @(18, @(1.02, 1.078252291480623, 1.072250280000458, 1.079906367714082, 1.086847122599561, 1.053510108015427, 1.084087873692681, 1.075454169336086, 1.083086031179812, 1.09000514589319, 1.085627402819942)),
    @(19, @(1.02, 1.07844593256937, 1.072348302957273, 1.080068479769913, 1.086987606004737, 1.053548123791993, 1.084254711937165, 1.075538122314645, 1.083233987846188, 1.090131645785223, 1.08579447799387)),
    @(20, @(1.02, 1.077579337651438, 1.071909554998336, 1.079342944218594, 1.086358857167625, 1.05337780330512, 1.083507991774764, 1.075162255720002, 1.082571730948912, 1.08956540127608, 1.085046697403188)),
    @(21, @(1.02, 1.074751447693752, 1.070476628638221, 1.076974642422129, 1.084306275844689, 1.052818678054585, 1.081069962483881, 1.073933012932994, 1.080408678780427, 1.08771545066787, 1.082605205831066)),
    @(22, @(1.02, 1.072965128186286, 1.069570593355378, 1.075478093378174, 1.083009085204294, 1.052462944782016, 1.079528903552637, 1.073154470722206, 1.079040824726038, 1.086545227851376, 1.081061958419377)),
    @(23, @(1.02, 1.073912894889243, 1.070051389813801, 1.076272167318139, 1.08369739195482, 1.052651926210303, 1.080346637438122, 1.073567734550254, 1.079766707829306, 1.087166267183206, 1.081880853580775)),
    @(24, @(1.02, 1.077626610437636, 1.071933493292593, 1.079382524943851, 1.086393158643216, 1.053387106857382, 1.083548730353714, 1.075182769498525, 1.082607864486626, 1.089596298144722, 1.085087493835593)),
    @(25, @(1.02, 1.081901651799289, 1.074096046230528, 1.082960604043421, 1.089493569387197, 1.054222339969202, 1.087230442640307, 1.077032893521622, 1.08587194140761, 1.09238639492479, 1.088774434575931))
)

foreach ($entry in $data) {
    $r = $entry[0]
    $vals = $entry[1]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value = $vals[$i]
    }
}

Write-Output "Updated $($data.Length) rows."